$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 768, shifting rows 768:809 down to 769:810
$ws.Rows.Item(768).Insert()

# Populate the newly inserted row 768 with the new data point.
# Force column A to Text format first so the date-like string "2026/02/07"
# is stored verbatim instead of being auto-converted to a date serial.
$ws.Range("A768").NumberFormat = "@"
$ws.Range("A768").Value = "2026/02/07"
$ws.Range("B768").Value = "土"
$ws.Range("C768").Value = 8
$ws.Range("D768").Value = 77
